# Refresh cached Market Board pricing/profit snapshot values across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR leve-profit worksheets.
# (Mirrors the scheduled-runner data refresh described in the commit.)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2097351
$ws.Range("J17").Value = 2159030
$ws.Range("L17").Value = 6477090
$ws.Range("N17").Value = -6477426
$ws.Range("H28").Value = 851
$ws.Range("I28").Value = 886.9167
$ws.Range("K28").Value = 886.9167
$ws.Range("M28").Value = -401.9167
$ws.Range("H40").Value = 3253.4285
$ws.Range("J40").Value = 3896
$ws.Range("L40").Value = 3896
$ws.Range("N40").Value = -4246
$ws.Range("H135").Value = 7206.4
$ws.Range("I135").Value = 7206.4
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 64857.6
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -62322.6
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 666
$ws.Range("I5").Value = 674.1
$ws.Range("J5").Value = 652.5
$ws.Range("K5").Value = 674.1
$ws.Range("L5").Value = 652.5
$ws.Range("M5").Value = -562.1
$ws.Range("N5").Value = -876.5
$ws.Range("H32").Value = 162601.72
$ws.Range("I32").Value = 174137.9
$ws.Range("K32").Value = 174137.9
$ws.Range("M32").Value = -173850.9
$ws.Range("H61").Value = 8044.75
$ws.Range("I61").Value = 8090.25
$ws.Range("K61").Value = 8090.25
$ws.Range("M61").Value = -7878.25
$ws.Range("H122").Value = 3357.3
$ws.Range("I122").Value = 3175.2222
$ws.Range("K122").Value = 9525.6666
$ws.Range("M122").Value = -7075.6666
$ws.Range("H130").Value = 50900
$ws.Range("J130").Value = 50900
$ws.Range("L130").Value = 50900
$ws.Range("N130").Value = -60940
$ws.Range("H136").Value = 8044.75
$ws.Range("I136").Value = 8090.25
$ws.Range("K136").Value = 24270.75
$ws.Range("M136").Value = -21720.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 666
$ws.Range("I4").Value = 674.1
$ws.Range("J4").Value = 652.5
$ws.Range("K4").Value = 674.1
$ws.Range("L4").Value = 652.5
$ws.Range("M4").Value = -559.1
$ws.Range("N4").Value = -882.5
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H86").Value = 8552.666999999999
$ws.Range("I86").Value = 2666.6667
$ws.Range("K86").Value = 2666.6667
$ws.Range("M86").Value = -1543.6667
$ws.Range("H89").Value = 8552.666999999999
$ws.Range("I89").Value = 2666.6667
$ws.Range("K89").Value = 13333.3335
$ws.Range("M89").Value = -7717.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 116.333336
$ws.Range("I7").Value = 147.25
$ws.Range("K7").Value = 147.25
$ws.Range("M7").Value = -34.25
$ws.Range("H94").Value = 1177.1111
$ws.Range("I94").Value = 879
$ws.Range("K94").Value = 879
$ws.Range("M94").Value = -428
$ws.Range("H132").Value = 13336544
$ws.Range("I132").Value = 2999.75
$ws.Range("K132").Value = 8999.25
$ws.Range("M132").Value = -6469.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1759.0588
$ws.Range("I12").Value = 1099
$ws.Range("J12").Value = 1847.0667
$ws.Range("K12").Value = 3297
$ws.Range("L12").Value = 5541.2001
$ws.Range("M12").Value = -3124
$ws.Range("N12").Value = -5887.2001
$ws.Range("H39").Value = 5613.263
$ws.Range("J39").Value = 5814
$ws.Range("L39").Value = 17442
$ws.Range("N39").Value = -18030
$ws.Range("H103").Value = 6541.3
$ws.Range("I103").Value = 108.333336
$ws.Range("J103").Value = 9298.286
$ws.Range("K103").Value = 325.000008
$ws.Range("L103").Value = 27894.858
$ws.Range("M103").Value = 553.999992
$ws.Range("N103").Value = -29652.858
$ws.Range("H107").Value = 41667000
$ws.Range("I107").Value = 896.5
$ws.Range("J107").Value = 50000220
$ws.Range("K107").Value = 2689.5
$ws.Range("L107").Value = 150000660
$ws.Range("M107").Value = -769.5
$ws.Range("N107").Value = -150004500
$ws.Range("H121").Value = 1321.6364
$ws.Range("I121").Value = 850.6667
$ws.Range("J121").Value = 1498.25
$ws.Range("K121").Value = 2552.0001
$ws.Range("L121").Value = 4494.75
$ws.Range("N121").Value = -7114.75
$ws.Range("M121").Value = -1242.0001
$ws.Range("H131").Value = 8123812.5
$ws.Range("I131").Value = 22730536
$ws.Range("J131").Value = 5341579.5
$ws.Range("K131").Value = 68191608
$ws.Range("L131").Value = 16024738.5
$ws.Range("M131").Value = -68186568
$ws.Range("N131").Value = -16034818.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 63200.26
$ws.Range("I80").Value = 88538.38
$ws.Range("K80").Value = 88538.38
$ws.Range("M80").Value = -87540.38
$ws.Range("H83").Value = 63200.26
$ws.Range("I83").Value = 88538.38
$ws.Range("K83").Value = 442691.9
$ws.Range("M83").Value = -437699.9
$ws.Range("H97").Value = 2150
$ws.Range("I97").Value = 2150
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2150
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1654
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 12565875
$ws.Range("I102").Value = 12852179
$ws.Range("K102").Value = 12852179
$ws.Range("M102").Value = -12850557
$ws.Range("H107").Value = 167831.33
$ws.Range("I107").Value = 250372.25
$ws.Range("K107").Value = 250372.25
$ws.Range("M107").Value = -248452.25
$ws.Range("H122").Value = 3656
$ws.Range("I122").Value = 1560.375
$ws.Range("K122").Value = 4681.125
$ws.Range("M122").Value = -2231.125
$ws.Range("H140").Value = 154990
$ws.Range("J140").Value = 154990
$ws.Range("L140").Value = 154990
$ws.Range("N140").Value = -165350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5884939.5
$ws.Range("I61").Value = 7694755.5
$ws.Range("K61").Value = 7694755.5
$ws.Range("M61").Value = -7694553.5
$ws.Range("H93").Value = 1543.4166
$ws.Range("J93").Value = 3500
$ws.Range("L93").Value = 3500
$ws.Range("N93").Value = -5996
$ws.Range("H94").Value = 58957.875
$ws.Range("J94").Value = 67499.5
$ws.Range("L94").Value = 67499.5
$ws.Range("N94").Value = -68851.5
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H113").Value = 5884939.5
$ws.Range("I113").Value = 7694755.5
$ws.Range("K113").Value = 7694755.5
$ws.Range("M113").Value = -7692585.5
$ws.Range("H122").Value = 5910.609
$ws.Range("J122").Value = 6312.8423
$ws.Range("L122").Value = 18938.5269
$ws.Range("N122").Value = -23838.5269
$ws.Range("H136").Value = 4200.619
$ws.Range("I136").Value = 3766.0833
$ws.Range("K136").Value = 11298.2499
$ws.Range("M136").Value = -8748.249899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 27649.334
$ws.Range("I122").Value = 32949.4
$ws.Range("K122").Value = 98848.20000000001
$ws.Range("M122").Value = -96398.20000000001
$ws.Range("H125").Value = 67723.17999999999
$ws.Range("J125").Value = 67723.17999999999
$ws.Range("L125").Value = 67723.17999999999
$ws.Range("N125").Value = -77563.17999999999
$ws.Range("H136").Value = 25314.28
$ws.Range("I136").Value = 49059.953
$ws.Range("K136").Value = 147179.859
$ws.Range("M136").Value = -144629.859
$ws.Range("H137").Value = 81663.336
$ws.Range("J137").Value = 81663.336
$ws.Range("L137").Value = 81663.336
$ws.Range("N137").Value = -91863.336
$ws.Range("H141").Value = 226901.33
$ws.Range("J141").Value = 226901.33
$ws.Range("L141").Value = 226901.33
$ws.Range("N141").Value = -237261.33
